$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.3425203228057133
$ws.Range("C2").Value = 0.1370775318625127
$ws.Range("D2").Value = 1.367787377993229

$ws.Range("B3").Value = -0.0860618221422374
$ws.Range("C3").Value = 1.383680499826445

$ws.Range("B4").Value = 1.32438096272888

$ws.Range("B5").Value = 1.231661280012998
$ws.Range("C5").Value = 0.1923542655531081
$ws.Range("D5").Value = 0.2120950120634008
$ws.Range("E5").Value = 0.3792844388692188

$ws.Range("B6").Value = 0.4395685548510502
$ws.Range("C6").Value = 0.3295146744469067
$ws.Range("D6").Value = 0.2496258682164595

$ws.Range("B7").Value = 0.4940016924669799
$ws.Range("C7").Value = 0.2435019605816055

$ws.Range("B8").Value = 0.3242193037695071

$ws.Range("B9").Value = 0.394042449657095
$ws.Range("C9").Value = 0.3427404628639549
$ws.Range("D9").Value = 0.2242390996078211
$ws.Range("E9").Value = 0.055485660899395

$ws.Range("B10").Value = 0.5010464375566571
$ws.Range("C10").Value = 0.2455732575174918
$ws.Range("D10").Value = 0.0184750902009912

$ws.Range("B11").Value = 0.4070253497240054
$ws.Range("C11").Value = 0.0588001744469144

$ws.Range("B12").Value = 0.1842008206034934

$ws.Range("B13").Value = -0.3665930774731743
$ws.Range("C13").Value = 0.101671561665663
$ws.Range("D13").Value = 0.04785045814007061

$ws.Range("B14").Value = -0.1157322460162614
$ws.Range("C14").Value = 0.09802926010601931

$ws.Range("B15").Value = 0.2279386281717039
